$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Activate()

$xlPasteFormats = [Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats

# --- Swap rows 19 and 20 (Mole <-> Mice), including row height ---
$a19 = $ws.Range("A19").Value2
$d19 = $ws.Range("D19").Value2

$a20 = $ws.Range("A20").Value2
$d20 = $ws.Range("D20").Value2

$ws.Range("A19").Value = $a20
$ws.Range("D19").Value = $d20
$ws.Rows.Item(19).RowHeight = 17.25

$ws.Range("A20").Value = $a19
$ws.Range("D20").Value = $d19
$ws.Rows.Item(20).AutoFit()

# --- Swap rows 21 and 22 (boss grasshopper <-> boss fly), including A-cell style/alignment ---
$a21 = $ws.Range("A21").Value2
$c21 = $ws.Range("C21").Value2
$d21 = $ws.Range("D21").Value2
$e21 = $ws.Range("E21").Value2

$a22 = $ws.Range("A22").Value2
$c22 = $ws.Range("C22").Value2
$d22 = $ws.Range("D22").Value2
$e22 = $ws.Range("E22").Value2

# Copy cell formats (keeps style indices clean instead of spawning new ones)
$ws.Range("A21").Copy()
$ws.Range("Z1").PasteSpecial($xlPasteFormats)   # stash A21's (boss grasshopper) format temporarily
$ws.Range("A22").Copy()
$ws.Range("A21").PasteSpecial($xlPasteFormats)  # A21 gets A22's (boss fly) format
$ws.Range("Z1").Copy()
$ws.Range("A22").PasteSpecial($xlPasteFormats)  # A22 gets the stashed boss grasshopper format
$ws.Range("Z1").Clear()

$ws.Range("A21").Value = $a22
$ws.Range("C21").Value = $c22
$ws.Range("D21").Value = $d22
$ws.Range("E21").Value = $e22

$ws.Range("A22").Value = $a21
$ws.Range("C22").Value = $c21
$ws.Range("D22").Value = $d21
$ws.Range("E22").Value = $e21

$excel.CutCopyMode = $false

# --- Update view state: scroll position and selection ---
$excel.ActiveWindow.ScrollRow = 7
$excel.ActiveWindow.ScrollColumn = 1
$ws.Rows.Item(26).Select()
